$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "26.416.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.696.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  +0.40%  "
$ws.Range("E4").Value2 = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "219.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.5495"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +4.37%  "
$ws.Range("E7").Value2 = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2755"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +1.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06467"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  +0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "22.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -0.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07688"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +2.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "1.696.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "4.549"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.5851"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.000008357"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -2.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "65.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "26.477.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "4.942"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "11.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "192.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "6.261"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +0.65%  "
$ws.Range("E23").Value2 = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "148.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +2.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.1331"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +7.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "7.928"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +2.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "15.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.06305"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -5.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "1.383"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +2.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.332"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "3.611"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.612"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.690"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +1.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.045"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.6166"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -1.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "2.412"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "2.721"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.01657"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +2.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "6.187"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  -2.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.118.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.8842"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -0.28%  "
$ws.Range("E42").Value2 = "  -0.20%  "
$ws.Range("E43").Value2 = "  +0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "1.844.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "57.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +1.38%  "
$ws.Range("E46").Value2 = "  -6.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "8.250"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +1.24%  "
$ws.Range("E48").Value2 = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.05277"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "6.120"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.4302"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -0.04%  "
